$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3206.0625
$ws.Range("I40").Value = 2764.4443
$ws.Range("J40").Value = 3773.8572
$ws.Range("K40").Value = 2764.4443
$ws.Range("L40").Value = 3773.8572
$ws.Range("M40").Value = -2589.4443
$ws.Range("N40").Value = -4123.8572

$ws.Range("H49").Value = 1083.2
$ws.Range("I49").Value = 708.5
$ws.Range("J49").Value = 1333
$ws.Range("K49").Value = 2125.5
$ws.Range("L49").Value = 3999
$ws.Range("M49").Value = -1989.5
$ws.Range("N49").Value = -4271

$ws.Range("H64").Value = 3270.9375
$ws.Range("I64").Value = 3062.75
$ws.Range("J64").Value = 3479.125
$ws.Range("K64").Value = 3062.75
$ws.Range("L64").Value = 3479.125
$ws.Range("M64").Value = -2814.75
$ws.Range("N64").Value = -3975.125

$ws.Range("H67").Value = 3270.9375
$ws.Range("I67").Value = 3062.75
$ws.Range("J67").Value = 3479.125
$ws.Range("K67").Value = 3062.75
$ws.Range("L67").Value = 3479.125
$ws.Range("M67").Value = -2204.75
$ws.Range("N67").Value = -5195.125

$ws.Range("H74").Value = 3875.08
$ws.Range("I74").Value = 3778.1304
$ws.Range("J74").Value = 4990
$ws.Range("K74").Value = 3778.1304
$ws.Range("L74").Value = 4990
$ws.Range("M74").Value = -2842.1304
$ws.Range("N74").Value = -6862

$ws.Range("H76").Value = 2897.8857
$ws.Range("I76").Value = 2756.7778
$ws.Range("J76").Value = 3374.125
$ws.Range("K76").Value = 2756.7778
$ws.Range("L76").Value = 3374.125
$ws.Range("M76").Value = -2441.7778
$ws.Range("N76").Value = -4004.125

$ws.Range("H77").Value = 3875.08
$ws.Range("I77").Value = 3778.1304
$ws.Range("J77").Value = 4990
$ws.Range("K77").Value = 18890.652
$ws.Range("L77").Value = 24950
$ws.Range("M77").Value = -14210.652
$ws.Range("N77").Value = -34310

$ws.Range("H79").Value = 2897.8857
$ws.Range("I79").Value = 2756.7778
$ws.Range("J79").Value = 3374.125
$ws.Range("K79").Value = 2756.7778
$ws.Range("L79").Value = 3374.125
$ws.Range("M79").Value = -1664.7778
$ws.Range("N79").Value = -5558.125

$ws.Range("H113").Value = 3436.4092
$ws.Range("I113").Value = 2470.3845
$ws.Range("J113").Value = 4831.778
$ws.Range("K113").Value = 2470.3845
$ws.Range("L113").Value = 4831.778
$ws.Range("M113").Value = 783.6154999999999
$ws.Range("N113").Value = -11339.778

$ws.Range("H137").Value = 1630.1842
$ws.Range("J137").Value = 986.6087
$ws.Range("L137").Value = 2959.8261
$ws.Range("N137").Value = -8059.8261

$ws.Range("H139").Value = 52995
$ws.Range("J139").Value = 52995
$ws.Range("L139").Value = 52995
$ws.Range("N139").Value = -63275

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H69").Value = 115000
$ws.Range("J69").Value = 180000
$ws.Range("L69").Value = 180000
$ws.Range("N69").Value = -181498

$ws.Range("H72").Value = 115000
$ws.Range("J72").Value = 180000
$ws.Range("L72").Value = 540000
$ws.Range("N72").Value = -547488

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1047.9412
$ws.Range("I22").Value = 1107.0625
$ws.Range("J22").Value = 102
$ws.Range("K22").Value = 1107.0625
$ws.Range("L22").Value = 102
$ws.Range("M22").Value = -934.0625
$ws.Range("N22").Value = -448

$ws.Range("H94").Value = 764.8570999999999
$ws.Range("I94").Value = 652.9474
$ws.Range("J94").Value = 1001.1111
$ws.Range("K94").Value = 652.9474
$ws.Range("L94").Value = 1001.1111
$ws.Range("M94").Value = -201.9474
$ws.Range("N94").Value = -1903.1111

$ws.Range("H105").Value = 2712.6667
$ws.Range("I105").Value = 2000
$ws.Range("J105").Value = 2855.2
$ws.Range("K105").Value = 2000
$ws.Range("L105").Value = 2855.2
$ws.Range("M105").Value = -253
$ws.Range("N105").Value = -6349.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1885.1548
$ws.Range("I31").Value = 1152
$ws.Range("J31").Value = 3615.4
$ws.Range("K31").Value = 1152
$ws.Range("L31").Value = 3615.4
$ws.Range("M31").Value = -857
$ws.Range("N31").Value = -4205.4

$ws.Range("H34").Value = 1885.1548
$ws.Range("I34").Value = 1152
$ws.Range("J34").Value = 3615.4
$ws.Range("K34").Value = 1152
$ws.Range("L34").Value = 3615.4
$ws.Range("M34").Value = -950
$ws.Range("N34").Value = -4019.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 920.2941
$ws.Range("I5").Value = 766.9231
$ws.Range("J5").Value = 1418.75
$ws.Range("K5").Value = 2300.7693
$ws.Range("L5").Value = 4256.25
$ws.Range("M5").Value = -2188.7693
$ws.Range("N5").Value = -4480.25

$ws.Range("H115").Value = 3942.6667
$ws.Range("I115").Value = 3942.6667
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 11828.0001
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = -10653.0001
$ws.Range("N115").Value = $null

$ws.Range("H122").Value = 3043.122
$ws.Range("I122").Value = 1038
$ws.Range("J122").Value = 3321.611
$ws.Range("K122").Value = 9342
$ws.Range("L122").Value = 29894.499
$ws.Range("M122").Value = -6892
$ws.Range("N122").Value = -34794.499

$ws.Range("H135").Value = 920.2941
$ws.Range("I135").Value = 766.9231
$ws.Range("J135").Value = 1418.75
$ws.Range("K135").Value = 6902.3079
$ws.Range("L135").Value = 12768.75
$ws.Range("M135").Value = -4367.3079
$ws.Range("N135").Value = -17838.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7634.1816
$ws.Range("I70").Value = 9096.857
$ws.Range("J70").Value = 5074.5
$ws.Range("K70").Value = 9096.857
$ws.Range("L70").Value = 5074.5
$ws.Range("M70").Value = -8826.857
$ws.Range("N70").Value = -5614.5

$ws.Range("H73").Value = 7634.1816
$ws.Range("I73").Value = 9096.857
$ws.Range("J73").Value = 5074.5
$ws.Range("K73").Value = 9096.857
$ws.Range("L73").Value = 5074.5
$ws.Range("M73").Value = -8160.857
$ws.Range("N73").Value = -6946.5

$ws.Range("H80").Value = 2501.3
$ws.Range("I80").Value = 2420.2354
$ws.Range("J80").Value = 2607.3076
$ws.Range("K80").Value = 2420.2354
$ws.Range("L80").Value = 2607.3076
$ws.Range("M80").Value = -1422.2354
$ws.Range("N80").Value = -4603.3076

$ws.Range("H83").Value = 2501.3
$ws.Range("I83").Value = 2420.2354
$ws.Range("J83").Value = 2607.3076
$ws.Range("K83").Value = 12101.177
$ws.Range("L83").Value = 13036.538
$ws.Range("M83").Value = -7109.177
$ws.Range("N83").Value = -23020.538

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 294.45456
$ws.Range("I55").Value = 242.375
$ws.Range("K55").Value = 242.375
$ws.Range("M55").Value = -69.375
$ws.Range("N55").Value = -779.33334

Write-Host "Updated Mandragora Profits values across sheets."